# Applies the cryptos.xlsx data refresh described in the commit.
# Updates Price (D) and Volume(1h) (E) columns for the latest rows,
# including the Uniswap/Polygon row swap (rows 22-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.758.58"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.546.02"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'607.46"
$ws.Range("D6").Value = "'174.34"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("D8").Value = "3.541.42"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D10").Value = "'0.201"
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "'0.587"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("D13").Value = "'47.68"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'0.0000281"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "4.115.12"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "'629.96"
$ws.Range("E16").Value = "  -6.95%  "
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "70.843.06"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "3.537.38"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.891"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  -10.57%  "
$ws.Range("D24").Value = "'15.98"
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("D25").Value = "'97.12"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'3.86"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").Value = "'2.62"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "'9.26"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'33.47"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'3.14"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'8.48"
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("D33").Value = "'1.34"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").Value = "'7.06"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").Value = "'568.17"
$ws.Range("E35").Value = "  -4.69%  "
$ws.Range("D36").Value = "'3.65"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "'10.81"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.102"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "'57.55"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("D44").Value = "3.343.78"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("E45").Value = "  +4.32%  "
$ws.Range("D46").Value = "0.0₃0720"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").Value = "'33.20"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").Value = "'2.66"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").Value = "'134.32"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "'5.72"
$ws.Range("E51").Value = "  -2.02%  "
